# The worksheet contained several placeholder "nan" rows (rows 21-27)
# followed by a TOTAL summary row (row 28). Remove the placeholder rows
# so the TOTAL row becomes row 21, shrinking the used range from
# A1:L28 down to A1:L21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("21:27").Delete() | Out-Null
